$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Insert two new rows right after the existing worker row (16), shifting the
# footer rows (21/22 -> 23/24) and their merged cells down automatically.
$ws.Range("A17:A18").EntireRow.Insert(-4121)

# Copy the formatting (styles/borders) of the existing data row into the two
# newly inserted rows so they match the look of the table.
$ws.Range("B16:J16").Copy() | Out-Null
$ws.Range("B17:J18").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# Fill in the two new worker records.
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "30317588"
$ws.Range("D17").Value = "CLAUDIA PATRICIA DUQUE RESTREPO"
$ws.Range("E17").Value = "2508"
$ws.Range("F17").Value = 755227
$ws.Range("G17").Value = 18880676

$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1047475286"
$ws.Range("D18").Value = "SANDRA MILENA CAMPO YUNES"
$ws.Range("E18").Value = "2508"
$ws.Range("F18").Value = 86957
$ws.Range("G18").Value = 2173913

# Update the summary figures: total overdue amount, worker count, period count.
$ws.Range("E11").Value = 870184
$ws.Range("C13").Value = 3
$ws.Range("F13").Value = 2

$wb.Save()
